# Update the two-digit-divided-by-one-digit practice table.
# The underlying row/cell grid shape is unchanged (20 rows x 5 cols); only
# the division problems shown in the 5 "data" rows (1, 5, 9, 13, 17) are
# refreshed. Addressing cells by (row, col) avoids any ambiguity from
# Find/Replace when a new value happens to equal another cell's old value
# (e.g. "70÷5=" is both an old value in row 17 and a new value in row 1).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$values = @{
    1  = @("42÷3=", "41÷4=", "70÷5=", "11÷5=", "98÷4=")
    5  = @("90÷5=", "72÷6=", "20÷8=", "55÷2=", "43÷8=")
    9  = @("15÷8=", "51÷2=", "28÷8=", "78÷9=", "96÷9=")
    13 = @("44÷8=", "80÷7=", "35÷2=", "73÷7=", "19÷3=")
    17 = @("25÷7=", "73÷6=", "17÷8=", "60÷4=", "65÷5=")
}

foreach ($rowIndex in $values.Keys) {
    $rowValues = $values[$rowIndex]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $rowValues[$col - 1]
    }
}
